$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shp = $master.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
Write-Output "Runs count: $($tr.Runs().Count)"
for ($i=1; $i -le $tr.Runs().Count; $i++) {
    $run = $tr.Runs($i)
    Write-Output "run $i : $($run.Text)"
}
